$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("D2").Value = "中山北路3300号Global HarborF5 世界港口小镇"
$ws1.Range("F2").Value = 580
$ws1.Range("F4").Value = 6366
$ws1.Range("F5").Value = 717
$ws1.Range("F6").Value = 1085
$ws1.Range("F7").Value = 70
$ws1.Range("F8").Value = 316
$ws1.Range("F10").Value = 8
$ws1.Range("F11").Value = 688
$ws1.Range("F12").Value = 1165
$ws1.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202402/1SB1YDwy1709194691582.jpeg"
$ws1.Range("F14").Value = 414
$ws1.Range("F17").Value = 1416
$ws1.Range("F18").Value = 665
$ws1.Range("F19").Value = 379
$ws1.Range("F20").Value = 396
$ws1.Range("F22").Value = 1066
$ws1.Range("F23").Value = 135
$ws1.Range("F24").Value = 2204
$ws1.Range("F25").Value = 250
$ws1.Range("F26").Value = 92
$ws1.Range("F27").Value = 393
$ws1.Range("F29").Value = 3559
$ws1.Range("F31").Value = 630

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value = 174
$ws2.Range("F8").Value = 705
$ws2.Range("F18").Value = 375
$ws2.Range("F20").Value = 4087
$ws2.Range("F24").Value = 189
$ws2.Range("F25").Value = 230
$ws2.Range("F32").Value = 1645

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F5").Value = 1187
$ws3.Range("F7").Value = 1573
$ws3.Range("F8").Value = 430
$ws3.Range("C10").Value = "上海·飘起来吧魔法泡泡-魔术表演（取消）"
$ws3.Range("F11").Value = 765
$ws3.Range("G10").Value = "不可售"

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 1187
$ws4.Range("F5").Value = 1573
$ws4.Range("F6").Value = 430
$ws4.Range("F8").Value = 765
$ws4.Range("D9").Value = "中山北路3300号Global HarborF5 世界港口小镇"
$ws4.Range("F9").Value = 580
$ws4.Range("F11").Value = 6366
$ws4.Range("F13").Value = 717
$ws4.Range("F14").Value = 1085
$ws4.Range("F15").Value = 705
$ws4.Range("F16").Value = 70
$ws4.Range("F17").Value = 316
$ws4.Range("F19").Value = 688
$ws4.Range("F22").Value = 1165
$ws4.Range("I22").Value = "//i0.hdslb.com/bfs/openplatform/202402/1SB1YDwy1709194691582.jpeg"
$ws4.Range("F24").Value = 414
$ws4.Range("F26").Value = 375
$ws4.Range("F28").Value = 1416
$ws4.Range("F30").Value = 665
$ws4.Range("F31").Value = 379
$ws4.Range("F32").Value = 396
$ws4.Range("F34").Value = 189
$ws4.Range("F35").Value = 230
$ws4.Range("F36").Value = 1066
$ws4.Range("F37").Value = 135
$ws4.Range("F39").Value = 2204
$ws4.Range("F41").Value = 1645
$ws4.Range("F42").Value = 250
$ws4.Range("F43").Value = 92
$ws4.Range("F44").Value = 393
$ws4.Range("F46").Value = 3559
$ws4.Range("F50").Value = 630
